$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "604.08")
# are preserved exactly as typed rather than being parsed into floating point
# numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.431.61"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "2.659.28"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "604.08"
$ws.Range("E5").Value = "  +1.47%  "

$ws.Range("D6").Value = "157.60"
$ws.Range("E6").Value = "  +4.43%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("E9").Value = "  +7.89%  "

$ws.Range("D10").Value = "0.408"
$ws.Range("E10").Value = "  +3.55%  "

$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("E12").Value = "  +1.53%  "

$ws.Range("D13").Value = "29.60"
$ws.Range("E13").Value = "  +5.86%  "

$ws.Range("E14").Value = "  +15.48%  "

$ws.Range("D15").Value = "3.138.55"
$ws.Range("E15").Value = "  +1.53%  "

$ws.Range("D16").Value = "65.204.26"
$ws.Range("E16").Value = "  +2.37%  "

$ws.Range("D17").Value = "2.658.74"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "12.80"
$ws.Range("E18").Value = "  +4.61%  "

$ws.Range("D19").Value = "4.91"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").Value = "360.23"
$ws.Range("E20").Value = "  +3.42%  "

$ws.Range("D21").Value = "7.36"
$ws.Range("E21").Value = "  +5.03%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "68.91"
$ws.Range("E23").Value = "  +2.12%  "

$ws.Range("D24").Value = "1.71"
$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("D25").Value = "9.53"
$ws.Range("E25").Value = "  +3.57%  "

$ws.Range("E26").Value = "  +16.58%  "

$ws.Range("D27").Value = "1.65"
$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("E29").Value = "  +2.07%  "

$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  +7.17%  "

$ws.Range("D31").Value = "540.24"
$ws.Range("E31").Value = "  -2.17%  "

$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +4.75%  "

$ws.Range("D35").Value = "6.43"
$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("E36").Value = "  +3.96%  "

$ws.Range("D37").Value = "20.62"
$ws.Range("E37").Value = "  +4.18%  "

$ws.Range("D38").Value = "163.01"
$ws.Range("E38").Value = "  -0.84%  "

$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "42.42"
$ws.Range("E42").Value = "  +6.10%  "

$ws.Range("D43").Value = "166.55"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").Value = "4.19"
$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("E45").Value = "  +7.13%  "

$ws.Range("D46").Value = "0.0617"
$ws.Range("E46").Value = "  +5.39%  "

$ws.Range("D47").Value = "23.12"
$ws.Range("E47").Value = "  -2.23%  "

$ws.Range("D48").Value = "0.659"
$ws.Range("E48").Value = "  +3.26%  "

$ws.Range("E49").Value = "  +4.64%  "

$ws.Range("E50").Value = "  +1.93%  "

$ws.Range("E51").Value = "  +3.26%  "

# Restore the default (General) style so the cells keep their original
# appearance/style index; only their text content has changed.
$ws.Range("D2:D51").Style = "Normal"
